$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6
# from serial date 45243 (2023-11-13) to 45244 (2023-11-14)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45244
}
